# Generate Report for Handback
# Adds a 3rd handback-status row (for 3074847d-c8c4-4e29-b8c5-dc8af17a4364)
# to the Overview / zh-cn / de-de tables.

$wb = $excel.ActiveWorkbook

$uuid = "3074847d-c8c4-4e29-b8c5-dc8af17a4364"

# Hyperlink targets, following the exact pattern already used by the
# existing rows (source repo commit is shared across all three sheets,
# target-repo commits are per-language).
$srcUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8817dadca1059cb3dd25a197cb0040db1be0431f/e2e/$uuid.md"
$zhcnUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/7f8003f2934ec718fc832705e40e1b1be4d15599/e2e/$uuid.md"
$dedeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/0947c4be144a63268fceaf2980ff09e957a3af67/e2e/$uuid.md"

$statusInSync = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet (sheet1) — new row 4
# Columns: A File Name | B Path And Name | C Extension | D Publish URL
#          E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "$uuid.md"
$wsOverview.Range("B4").Value = "e2e\$uuid.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = "2016-08-27 04:43:07"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcUrl, $null, $null, "e2e\$uuid.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet (sheet2) — new row 4
# Columns: A Source File Name | B File Extension | C Status | D Source Path
#          E Priority | F Content Duplicate | G Correspond Handoff File
#          H Correspond Handoff Datetime | I Target File
#          J Correspond Handback File | K Correspond Handback DateTime
#          L Reference Tokens | M To be localized | N Dependency From
#          O Has metadata | P Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "$uuid.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusInSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "$uuid.c278ded7a50109a917f550759164a20555dd2375.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-27 04:42:59"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = "$uuid.md"
$wsZhCn.Range("J4").Value = "$uuid.c278ded7a50109a917f550759164a20555dd2375.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-27 04:43:27"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $srcUrl, $null, $null, "$uuid.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), $zhcnUrl, $null, $null, "$uuid.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (sheet3) — new row 4 (same column layout as zh-cn)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "$uuid.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusInSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "$uuid.c278ded7a50109a917f550759164a20555dd2375.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-27 04:43:07"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = "$uuid.md"
$wsDeDe.Range("J4").Value = "$uuid.c278ded7a50109a917f550759164a20555dd2375.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-27 04:43:33"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $srcUrl, $null, $null, "$uuid.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), $dedeUrl, $null, $null, "$uuid.md") | Out-Null
